$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2026-02-09 Monday", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-10 Tuesday", 2)

# Update the division problems in the table (rows 1, 5, 9, 13, 17 hold the
# visible text; the other rows are blank spacer rows for handwriting).
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "30÷7="
$t.Cell(1,2).Range.Text  = "83÷3="
$t.Cell(1,3).Range.Text  = "61÷4="
$t.Cell(1,4).Range.Text  = "54÷4="
$t.Cell(1,5).Range.Text  = "89÷4="

$t.Cell(5,1).Range.Text  = "11÷9="
$t.Cell(5,2).Range.Text  = "83÷7="
$t.Cell(5,3).Range.Text  = "84÷7="
$t.Cell(5,4).Range.Text  = "64÷4="
$t.Cell(5,5).Range.Text  = "65÷7="

$t.Cell(9,1).Range.Text  = "34÷9="
$t.Cell(9,2).Range.Text  = "57÷9="
$t.Cell(9,3).Range.Text  = "18÷8="
$t.Cell(9,4).Range.Text  = "27÷8="
$t.Cell(9,5).Range.Text  = "89÷2="

$t.Cell(13,1).Range.Text = "17÷8="
$t.Cell(13,2).Range.Text = "90÷8="
$t.Cell(13,3).Range.Text = "60÷2="
$t.Cell(13,4).Range.Text = "55÷7="
$t.Cell(13,5).Range.Text = "26÷8="

$t.Cell(17,1).Range.Text = "12÷2="
$t.Cell(17,2).Range.Text = "11÷5="
$t.Cell(17,3).Range.Text = "17÷5="
$t.Cell(17,4).Range.Text = "55÷5="
$t.Cell(17,5).Range.Text = "42÷4="
